$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") values of 52 were dropped a leading "2" digit;
# recover the original values of 252 for all data rows (2-65).
$ws.Range("B2:B65").Value = 252
